$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''90.468.28'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.93%  '
$ws.Range('D3').Value = '''3.135.49'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.12%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').Value = '''215.26'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Value = '''621.53'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.40%  '
$ws.Range('D7').Value = '''1.13'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +26.64%  '
$ws.Range('D8').Value = '''0.362'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.84%  '
$ws.Range('D9').Value = '''1.00'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('D10').Value = '''3.133.93'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.15%  '
$ws.Range('D11').Value = '''0.745'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +7.52%  '
$ws.Range('E12').Value = '  +5.78%  '
$ws.Range('B13').Value = 'ShibaInu'
$ws.Range('C13').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D13').Value = '''0.0000244'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.27%  '
$ws.Range('B14').Value = 'Toncoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D14').Value = '''5.65'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.97%  '
$ws.Range('D15').Value = '''35.02'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.59%  '
$ws.Range('D16').Value = '''90.321.13'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.73%  '
$ws.Range('D17').Value = '''3.715.87'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.36%  '
$ws.Range('D18').Value = '''3.140.32'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.90%  '
$ws.Range('D19').Value = '''3.70'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.24%  '
$ws.Range('D20').Value = '''14.49'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.89%  '
$ws.Range('B21').Value = 'PEPE'
$ws.Range('C21').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D21').Value = '''0.0000213'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -6.55%  '
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').Value = '''462.87'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +7.13%  '
$ws.Range('D23').Value = '''9.04'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.82%  '
$ws.Range('D24').Value = '''5.30'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.78%  '
$ws.Range('B25').Value = 'NEARProtocol'
$ws.Range('C25').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D25').Value = '''5.90'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +5.13%  '
$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D26').Value = '''94.95'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +13.13%  '
$ws.Range('D27').Value = '''12.23'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.77%  '
$ws.Range('D28').Value = '''3.307.32'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.38%  '
$ws.Range('D29').Value = '''1.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.19%  '
$ws.Range('D30').Value = '''0.163'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.63%  '
$ws.Range('D31').Value = '''9.18'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.20%  '
$ws.Range('D32').Value = '''0.210'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +43.25%  '
$ws.Range('D33').Value = '''26.63'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +15.61%  '
$ws.Range('D34').Value = '''516.25'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.80%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').Value = '''0.146'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.75%  '
$ws.Range('B36').Value = 'PancakeSwap'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D36').Value = '''1.93'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.17%  '
$ws.Range('D37').Value = '''6.98'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.65%  '
$ws.Range('E38').Value = '  +2.18%  '
$ws.Range('E39').Value = '  -7.72%  '
$ws.Range('D40').Value = '''0.0911'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +26.97%  '
$ws.Range('D41').Value = '''22.22'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.46%  '
$ws.Range('B42').Value = 'PolygonEcosystemToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D42').Value = '''0.424'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +14.50%  '
$ws.Range('B43').Value = 'Binance-PegBSC-USD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D43').Value = '''0.754'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -24.59%  '
$ws.Range('D44').Value = '''1.00'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.02%  '
$ws.Range('D45').Value = '''1.97'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.98%  '
$ws.Range('B46').Value = 'USDe'
$ws.Range('C46').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D46').Value = '''1.00'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').Value = '''0.730'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +19.72%  '
$ws.Range('D48').Value = '''4.74'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +12.61%  '
$ws.Range('D49').Value = '''150.43'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +6.41%  '
$ws.Range('B50').Value = 'OKB'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D50').Value = '''45.36'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.82%  '
$ws.Range('B51').Value = 'ImmutableX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D51').Value = '''1.36'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +9.01%  '
